$d = $word.ActiveDocument

# 1. Merge the split "Petra (waving grinning): Hello! About time you woke up,
#    don't you think?" line into a single run (text itself is unchanged).
$d.Content.Find.Execute("Petra (waving grinning): Hello! About time you woke up, don't you think?", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Petra (waving grinning): Hello! About time you woke up, don't you think?", 2)

# 2. "Petra (surprise indignant):" -> "Petra (surprise surprise):"
$d.Content.Find.Execute("Petra (surprise indignant):", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Petra (surprise surprise):", 2)

# 3. "Petra (neutral raised_eyebrow): 100% sure?" -> "Petra (arms_crossed skeptical): 100% sure?"
$d.Content.Find.Execute("Petra (neutral raised_eyebrow): 100% sure?", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Petra (arms_crossed skeptical): 100% sure?", 2)

# 4. Merge the split "Petra (neutral disappointed): Alright, well..." line into
#    a single run (text itself is unchanged).
$d.Content.Find.Execute("Petra (neutral disappointed): Alright, well", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Petra (neutral disappointed): Alright, well", 2)

# 5. Remove the stray standalone "Teacher (smoking sigh):" paragraph that
#    precedes the "However, instead of the office..." paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Teacher (smoking sigh):`r") {
        $target = $p
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# 6. "Teacher: It's nice to be young, huh?" -> "Teacher: Must be nice to be young, huh?"
$d.Content.Find.Execute("Teacher: It's nice to be young, huh?", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Teacher: Must be nice to be young, huh?", 2)

# 7. "And with that, she walks off" -> "And with that she walks off"
$d.Content.Find.Execute("And with that, she walks off", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "And with that she walks off", 2)
